$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as plain text in the sheet (inline/shared
# strings), not numbers, so they must be written back as text too. Forcing the
# NumberFormat to "@" (Text) before the write stops Excel from silently
# re-interpreting the numeric-looking string as a float (which would mangle
# values like 245.54 into 245.53999999999999). Resetting the style back to
# "Normal" immediately after removes the temporary formatting again so the
# cell's style is left exactly as it was before the edit.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "245.54"
Set-TextValue "D3" "23.98"
Set-TextValue "D4" "5.360"
Set-TextValue "D5" "0.05817"
Set-TextValue "D6" "6.473"
Set-TextValue "D7" "3.363"
Set-TextValue "D8" "0.8095"
Set-TextValue "D9" "0.9184"
Set-TextValue "D10" "0.1405"
Set-TextValue "D11" "0.07383"
Set-TextValue "D12" "0.03142"
Set-TextValue "D13" "0.03073"
Set-TextValue "D14" "0.09364"
Set-TextValue "D15" "3.864"
Set-TextValue "D16" "0.001547"
Set-TextValue "D17" "0.04694"
Set-TextValue "D18" "0.0005982"
Set-TextValue "D19" "0.006119"
Set-TextValue "D20" "0.001247"
Set-TextValue "D21" "0.004690"
Set-TextValue "D22" "0.00008798"
Set-TextValue "D23" "3.592"
Set-TextValue "D25" "0.3183"
Set-TextValue "D40" "0.03837"
Set-TextValue "D41" "0.003067"
Set-TextValue "D42" "0.1066"
Set-TextValue "D43" "0.002749"
Set-TextValue "D44" "0.009028"
Set-TextValue "D45" "0.00005247"
Set-TextValue "D47" "0.6872"
Set-TextValue "D48" "0.001827"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.0002001"

# Column E values are non-numeric text already, so a direct assignment is safe.
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"
$ws.Range("E48").Value = "47BOLOBOLO"

Write-Host "Applied symbol list update"
